$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted at row 71, pushing all the
# subsequent rows (old 71-135) down by one (new 72-136). Insert a blank
# row before row 71 first so everything below shifts down.
$ws.Rows("71:71").Insert()

# Populate the newly inserted row 71 with the new record's data.
$ws.Range("A71").Value = 11
$ws.Range("B71").Value = "Vega Monumental Concepción"
$ws.Range("C71").Value = "Bíobío"
$ws.Range("D71").Value = 44729
$ws.Range("E71").Value = 8
$ws.Range("F71").Value = 100112043
$ws.Range("G71").Value = "Pepino ensalada"
$ws.Range("H71").Value = "Sin especificar"
$ws.Range("I71").Value = "Primera"
$ws.Range("J71").Value = 220
$ws.Range("K71").Value = 18000
$ws.Range("L71").Value = 19000
$ws.Range("M71").Value = 18455
$ws.Range("N71").Value = "`$/caja 60 unidades"
$ws.Range("O71").Value = "Región de Arica y Parinacota"
$ws.Range("P71").Value = 308
$ws.Range("Q71").Value = 60
$ws.Range("R71").Value = "Hortaliza"
